# rebase, fix unique constraints and test
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename activitycode headers to avoid unique-constraint collisions.
$ws.Range("G1").Value = "activitycode_code_set1"
$ws.Range("H1").Value = "activitycode_code_set2"
$ws.Range("I1").Value = "Activitycode_3"

# Normalise the "unit" column (V) style so every cell shares the same
# formatting (removes a now-unused duplicate font/style).
$ws.Range("V1:V6").Style = $ws.Range("U1").Style

# Update default column width and selection to match the fixed test fixture.
$ws.StandardWidth = 11.625
$ws.Range("G11").Select()
